# expense.xlsx update — "all" sheet: clear stale "chiya samosa" entry in
# row 5, add two new expense rows (66432 / 66433), and move the selection
# cursor to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# --- Row 5: remove the old "chiya samosa" / 100 entry -----------------
# Date (A5) stays as-is; only description + amount are cleared.
$ws.Range("B5:C5").ClearContents()

# --- Row 6: date only for now (text filled in after row 7, below, so the
# shared-string table ends up with the same insertion order as the
# original author's edit: "dinesh vinaju ..." before "Petrol, chiya
# samosa"). ----------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 66432

# --- Row 7: new entry, wrapped remarks style (like B3) -----------------
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 66433

$ws.Range("B3").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = "dinesh vinaju & prabin chiya, irika bus fare, Photocopy, apple, chocolate"

$ws.Range("C7").Formula = "=65+100+50+325+50"
$ws.Rows.Item(7).RowHeight = 28.8

# --- Row 6: fill in description + amount now -------------------------
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = "Petrol, chiya samosa"

$ws.Range("C6").Formula = "=150+100"

# --- Move the selection cursor to C3 ----------------------------------
$ws.Activate() | Out-Null
$ws.Range("C3").Select() | Out-Null
